# Add a new "Cameras View" worksheet, positioned right before "Tiles View",
# listing the GAME_CAMERA architecture entry (mirrors the layout of the
# existing "Prefabs View" / "UI View" sheets: IDS | Architecture ID | Prefab resource path).

$wb = $excel.ActiveWorkbook

$tilesView = $wb.Worksheets.Item("Tiles View")
$camerasView = $wb.Worksheets.Add($tilesView)
$camerasView.Name = "Cameras View"

$camerasView.Range("A1").Value = "IDS"
$camerasView.Range("B1").Value = "Architecture ID"
$camerasView.Range("C1").Value = "Prefab resource path"

$camerasView.Range("A2").Value = "GAME_CAMERA"
$camerasView.Range("B2").Value = "GameCamera"
$camerasView.Range("C2").Value = "Prefabs/Camera/GameCamera"
